$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sample")

# Row 6 (new): add Guard character row - set B6/C6 first so new shared strings
# CHAR_GUARD and Guard are appended to the shared string table before WEAPON_STING
$ws.Range("A6").Value = 10002
$ws.Range("B6").Value = "CHAR_GUARD"
$ws.Range("C6").Value = "Guard"

# Row 4 (Id 10000, Thief): change Default_Weapon from WEAPON_SLASH_PROJECTILE to WEAPON_STING
$ws.Range("D4").Value = "WEAPON_STING"

$ws.Range("D6").Value = "WEAPON_STING"
$ws.Range("E6").Value = "None"
$ws.Range("F6").Value = "None"
